$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: Reverse Linked List (Iteratively) complexity rating bumped 3* -> 4*
# (a new related problem, Has Cycle, was added below it)
$ws.Range("E17").Value = "4*"

# Row 22: new entry "Has Cycle" (LinkedListV1.py)
$ws.Range("A22").Value = "LC"
$ws.Range("B22").Value = "Has Cycle "
$ws.Range("C22").Value = "Easy"
$ws.Range("E22").Value = 1
$ws.Range("K22").Value = "O(n)"
$ws.Range("L22").Value = "O(1)"
$ws.Range("N22").Value = "Yes"

# Match formats of analogous existing cells (numeric Python column / Yes-No column)
$ws.Range("E19").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("N4").Copy() | Out-Null
$ws.Range("N22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row heights grew slightly (18.75 -> 19.5) on the header + data rows following the edit
$ws.Range("A1:A2").EntireRow.RowHeight = 19.5
$ws.Range("A5:A29").EntireRow.RowHeight = 19.5
